$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-09-21 (45190) to 2023-09-23 (45192) for every data row (rows 2-270).
$ws.Range("C2:C270").Value = 45192
